$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds daily covid data for Campogalliano in A1:D464.
# Extend it through 2022-01-05 (27 more days), appending rows 465-491:
#   A = date (serial), B = nuovi pos., C = somma mobile 7gg., D = per 100k ab.

# Copy A464 (incl. its date number format / style) onto the new date column
# cells so they inherit the same formatting as the rest of column A.
$ws.Range("A464").Copy($ws.Range("A465:A491"))

$newData = @(
    @(465, 44539, 4, 20, 229.9908003679853),
    @(466, 44540, 4, 21, 241.4903403863846),
    @(467, 44541, 0, 17, 195.4921803127875),
    @(468, 44542, 6, 22, 252.9898804047838),
    @(469, 44543, 8, 25, 287.4885004599816),
    @(470, 44544, 3, 26, 298.9880404783809),
    @(471, 44545, 0, 25, 287.4885004599816),
    @(472, 44546, 9, 30, 344.9862005519779),
    @(473, 44547, 4, 30, 344.9862005519779),
    @(474, 44548, 0, 30, 344.9862005519779),
    @(475, 44550, 9, 33, 379.4848206071757),
    @(476, 44551, 7, 32, 367.9852805887764),
    @(477, 44552, 0, 29, 333.4866605335786),
    @(478, 44553, 7, 36, 413.9834406623735),
    @(479, 44554, 1, 28, 321.9871205151794),
    @(480, 44555, 0, 24, 275.9889604415824),
    @(481, 44556, 5, 29, 333.4866605335786),
    @(482, 44557, 8, 28, 321.9871205151794),
    @(483, 44558, 10, 31, 356.4857405703772),
    @(484, 44559, 20, 51, 586.4765409383624),
    @(485, 44560, 10, 54, 620.9751609935603),
    @(486, 44561, 5, 58, 666.9733210671573),
    @(487, 44562, 3, 61, 701.4719411223551),
    @(488, 44563, 30, 86, 988.9604415823367),
    @(489, 44564, 29, 107, 1230.450781968721),
    @(490, 44565, 12, 109, 1253.44986200552),
    @(491, 44566, 11, 100, 1149.954001839926)
)

foreach ($r in $newData) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}

Write-Output "Appended $($newData.Count) rows (465-491) of covid data through 2022-01-05"
